$d = $word.ActiveDocument

# The "Requisitos" bullet paragraph lists three courses, one per line, each
# line being its own run (<w:r>) containing a <w:t> followed by a manual
# line break (<w:br/>):
#
#   LOB1021 -  Física IV  (Requisito)
#   LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
#   LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)
#
# The edit moves the "LOM3246" line from the end of the list to the very
# beginning, while "LOB1021" and "LOM3016" keep their original relative
# order (and, per the target XML, their original run objects are left
# completely untouched - only the LOM3246 run is removed from the end and
# re-inserted, as its own run, at the front).

$needle = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)"

$found = $d.Content
$found.Find.Text = $needle
$found.Find.Execute() | Out-Null

if (-not $found.Find.Found) {
    throw "Could not find the LOM3246 requirement line to move."
}

# Capture the text to move and the paragraph that currently contains it
# before any mutation invalidates the found range.
$lineText = $found.Text
$paraStart = $found.Paragraphs(1).Range.Start

# Remove the line together with its trailing manual line break (one char)
# from its current (last) position in the paragraph.
$toDelete = $d.Range($found.Start, $found.End + 1)
$toDelete.Delete()

# Re-insert it as a brand new run at the very start of the same paragraph
# (the other two lines' runs are never touched by this script, so they
# remain separate, unmerged runs, matching the target structure).
$insertion = $d.Range($paraStart, $paraStart)
$insertion.InsertBefore($lineText + [char]11)
